# "Aylık Plan" workbook refactor:
# The Başlangıç/Bitiş (start/end) columns (B, C) used to hold the dates as
# plain text strings ("01.10.2023", ...). Convert them to real Excel date
# serial values (with the existing m/d/yyyy date format) for the three
# task rows, which also drops the now-unused text strings from the shared
# string table. Finish by moving the active selection to C4, matching
# where the author's cursor ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Community Website Updates -> 01.10.2023 - 31.10.2023
$ws.Range("B2").Value = 45200
$ws.Range("C2").Value = 45230

# Row 3: News Craft -> 01.11.2023 - 30.11.2023
$ws.Range("B3").Value = 45231
$ws.Range("C3").Value = 45260

# Row 4: Social Master -> 01.12.2023 - 31.12.2023
$ws.Range("B4").Value = 45261
$ws.Range("C4").Value = 45291

# Keep the same date display format already used in column B/C.
$ws.Range("B2:C4").NumberFormat = "m/d/yyyy;@"

# Author's cursor ended up on C4 after the edit.
$ws.Range("C4").Select()
